# "Worked on temporal resolution"
# The Demand sheet ("t" / "EU27.Elec") previously modelled a single
# time step (t=1) carrying the full annual demand. It is extended here to
# twelve time steps (t=1..12), each carrying an equal share of that demand,
# so the model now resolves demand at a finer (monthly) temporal resolution.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand")

# Existing row 3 (t=1) keeps its value in sync with the new per-step amount.
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 421489583

# Add time steps 2 through 12, each with the same demand value.
$newRows = @(
  @(2, 421489583),
  @(3, 421489583),
  @(4, 421489583),
  @(5, 421489583),
  @(6, 421489583),
  @(7, 421489583),
  @(8, 421489583),
  @(9, 421489583),
  @(10, 421489583),
  @(11, 421489583),
  @(12, 421489583)
)

$row = 4
foreach ($pair in $newRows) {
  $ws.Cells.Item($row, 1).Value = $pair[0]
  $ws.Cells.Item($row, 2).Value = $pair[1]
  $row = $row + 1
}

# Column B now holds wider numbers - size it to fit (matches the
# "t"/value columns elsewhere in the workbook).
$ws.Columns.Item(2).ColumnWidth = 9.166666666666666

# Bring the Demand sheet to the front (it becomes the active/selected tab)
# and leave the cursor parked the way the author left it.
$ws.Activate() | Out-Null
$ws.Range("F12").Select() | Out-Null
